$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold + border + alignment) from the last existing header cell
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values for every data row (2-49)
$ws.Range("AD2:AD49").Value = 89
$ws.Range("AE2:AE49").Value = 73
$ws.Range("AF2:AF49").Value = 0

Write-Host "done"
